$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion note text (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.13 = 25026.84 pesos`n✅ 25026.84 pesos = 6.14 = 982.22 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet rates ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 163.025
$ws2.Range("O10").Value = 4080
$ws2.Range("O12").Value = 159.93
